$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D ("Tipo") to make room for "MAE"
$ws.Columns.Item(4).Insert()

# Copy header formatting (style) from C1 into new D1, then set its value
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("D1").Value = "MAE"

# New column D value (MAE metric)
$ws.Range("D2").Value = 0.3851510881690979

# Update existing B2 and C2 values (MSE, R2)
$ws.Range("B2").Value = 0.2394339242035501
$ws.Range("C2").Value = 0.9823097494285518
